# TC12_Canine_Filter_Breed-BelgMalin.xlsx — fill in the Neo4j query cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# The Cypher/Neo4j query text that goes into the (previously empty) B2 cell.
$query = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN ['Belgian Malinois'] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(s.clinical_study_designation,'') AS ``Study Code`` , coalesce(s.clinical_study_type,'') AS  ``Study Type``, coalesce(demo.breed,'') AS Breed , coalesce(diag.disease_term,'') AS Diagnosis , coalesce(diag.stage_of_disease,'') AS ``Stage of Disease`` ,  coalesce(demo.patient_age_at_enrollment,'') AS Age , coalesce(demo.sex,'') AS Sex , coalesce(demo.neutered_indicator,'') AS  ``Neutered Status``"

$ws.Range("B2").Value = $query

# Row 2 grows tall enough to show the (wrapped) query text.
$ws.Range("B2").EntireRow.RowHeight = 188.5

# The view now scrolls/selects column B instead of C.
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("B2").Select()
